# Update on 20250411 part 8
# Add the "江西教育" (Jiangxi Education) channel as two new rows (IPA + IPB)
# at the bottom of the multicast-source table, following the same
# channel/comma/protocol/ip-type/url layout used by every other row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48: 江西教育, ',', 'http://', 'IPA', /rtp/239.252.219.115:5140
$ws.Cells.Item(48, 1).Value = "江西教育"
$ws.Cells.Item(48, 2).Value = ","
$ws.Cells.Item(48, 3).Value = "http://"
$ws.Cells.Item(48, 4).Value = "IPA"
$ws.Cells.Item(48, 5).Value = "/rtp/239.252.219.115:5140"

# Row 49: 江西教育, ',', 'http://', 'IPB', /rtp/239.252.219.115:5140
$ws.Cells.Item(49, 1).Value = "江西教育"
$ws.Cells.Item(49, 2).Value = ","
$ws.Cells.Item(49, 3).Value = "http://"
$ws.Cells.Item(49, 4).Value = "IPB"
$ws.Cells.Item(49, 5).Value = "/rtp/239.252.219.115:5140"
